$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

for ($row = 2; $row -le 47; $row++) {
    $cell = $ws.Cells.Item($row, 4)
    $val = $cell.Value2
    if ($val -eq "MOLLY MCNINCH") {
        $cell.Value = "T"
    } elseif ($val -eq "STUDENT") {
        $cell.Value = "S"
    }
}
